# This document had its trial/payment-deadline dates advanced by one day
# in four places. Each occurrence lives in its own run with distinctive
# surrounding text, so we replace the most specific (longest/most unique)
# strings first, leaving the generic "June 11, 2022" (the bold, standalone
# occurrence) for last so it can't accidentally clobber the other matches.

$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2
$wdFindContinue = 1
$wdReplaceAll = 2

# 1. "Defendant appeared in Court on June 11, 2022" -> "...June 12, 2022"
$r1 = $d.Content.Find.Execute(
    "Defendant appeared in Court on June 11, 2022", $true, $false, $false,
    $false, $false, $true, $wdFindContinue, $false,
    "Defendant appeared in Court on June 12, 2022", $wdReplaceAll)

# 2. " license is suspended from June 11, 2022" -> "...June 12, 2022"
$r2 = $d.Content.Find.Execute(
    " license is suspended from June 11, 2022", $true, $false, $false,
    $false, $false, $true, $wdFindContinue, $false,
    " license is suspended from June 12, 2022", $wdReplaceAll)

# 3. The remaining standalone "June 11, 2022" (bold "Fines and Costs" date)
$r3 = $d.Content.Find.Execute(
    "June 11, 2022", $true, $false, $false,
    $false, $false, $true, $wdFindContinue, $false,
    "June 12, 2022", $wdReplaceAll)

# 4. "August 10, 2022" -> "August 11, 2022"
$r4 = $d.Content.Find.Execute(
    "August 10, 2022", $true, $false, $false,
    $false, $false, $true, $wdFindContinue, $false,
    "August 11, 2022", $wdReplaceAll)

Write-Output "Replace 1 (Defendant appeared...): $r1"
Write-Output "Replace 2 (license is suspended...): $r2"
Write-Output "Replace 3 (standalone June 11, 2022): $r3"
Write-Output "Replace 4 (August 10, 2022): $r4"
